$d = $word.ActiveDocument

# --- Fix 1: merge the split "(as ... a group) interview ..." sentence and its paragraph break ---

# Step 1: expand the run before the _GoBack bookmark to absorb the text that
# currently lives in the run right after the bookmark (fixing the stray
# double-space down to a single space in the process).
$old1 = "You as a group need to choose one of these two approaches. If it all possible, it is strongly recommended that you (as "
$new1 = "You as a group need to choose one of these two approaches. If it all possible, it is strongly recommended that you (as a group) interview an IT professional. Not only will this information be more direct, it may provide you with a contact that you will find helpful in the future. If you choose this option, you will be provided with a list of "
$null = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# Step 2: the text that used to sit in the run right after the bookmark is now
# duplicated (it was copied into the run above) -- clear that run out.
$old2 = "a group) interview an IT professional. Not only will this information be more direct, it may provide you with a contact that you will find helpful in the future. If you choose this option, you will be provided with a list of  "
$null = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Step 3: the paragraph that used to start with "questions to ask..." should
# become part of the same paragraph -- delete the paragraph mark joining them.
$p = $d.Paragraphs.Item(5)
$mark = $d.Range($p.Range.End - 1, $p.Range.End)
$mark.Delete()

# --- Fix 2: "Are they other IT professionals?" should be a single run, no proofErr wrapper ---
$old4 = "What kinds of people does the IT professional interact with? Are they other IT professionals? Clients? Investors? The general public? "
$null = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $old4, 2)
